$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

function Set-RowValues {
    param($ws, [int]$row, [object[]]$vals)
    $arr = New-Object 'object[,]' 1, $vals.Length
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $arr[0, $i] = $vals[$i]
    }
    $ws.Range($ws.Cells.Item($row, 4), $ws.Cells.Item($row, 13)).Value = $arr
}

# Insert two new columns before column D (new quarters), shifting existing
# D:K quarterly data right to F:M.
$ws.Range("D:E").Insert()

# The newly inserted D:E columns pick up the default style; carry over the
# number formats / styles that the (now shifted) column F cells use so the
# new quarter columns look like all the others (dates formatted as dates,
# financial figures as numbers).
$ws.Range("F5:F102").Copy()
$ws.Range("D5:E102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

Set-RowValues $ws 7 @(43465, 43373, 43281, 43190, 43100, 43008, 42916, 42825, 42735, 42643)
Set-RowValues $ws 8 @(2383000, 1819000, 1543000, 2242000, 1746000, 1403000, 1241000, 1211000, 571600, 465700)
Set-RowValues $ws 9 @(1543000, 1176000, 1017000, 1268000, 1042000, 949000, 810000, 741000, 321600, 324400)
Set-RowValues $ws 10 @(840000, 643000, 526000, 974000, 704000, 454000, 431000, 470000, 250000, 141300)
Set-RowValues $ws 12 @("NA", "NA", "NA", "NA", "NA", "NA", "NA", "NA", "NA", "NA")
Set-RowValues $ws 13 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws 14 @(0, 20000, 15000, 0, 4000, 34000, 33000, 48000, 65300, 52000)
Set-RowValues $ws 15 @(116000, 113000, 111000, 109000, 104000, 92000, 90000, 70000, 68000, 49200)
Set-RowValues $ws 17 @(1867000, 1406000, 1222000, 1495000, 1305000, 1131000, 1000000, 877000, 501800, 476200)
Set-RowValues $ws 18 @(516000, 413000, 321000, 747000, 441000, 272000, 241000, 334000, 69800, -10500)
Set-RowValues $ws 20 @(-71000, 26000, 32000, 74000, 51000, 2000, -31000, -6000, 237800, 29700)
Set-RowValues $ws 21 @(561000, 552000, 464000, 930000, 596000, 366000, 300000, 398000, 375600, 68400)
Set-RowValues $ws 22 @(210000, 209000, 206000, 206000, 208000, 186000, 188000, 156000, 158000, 148100)
Set-RowValues $ws 23 @(235000, 230000, 147000, 615000, 284000, 88000, 22000, 172000, 149500, -128800)
Set-RowValues $ws 24 @(12000, 3000, -3000, 15000, -400000, -2000, 1000, 0, 0, 1600)
Set-RowValues $ws 25 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws 26 @(223000, 227000, 150000, 600000, 684000, 90000, 21000, 172000, 149500, -130400)
Set-RowValues $ws 27 @(67000, 65000, -18000, 357000, 531000, -289000, -285000, 54000, 109700, -100400)
Set-RowValues $ws 28 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws 29 @("NA", "NA", "NA", "NA", -404000, "NA", "NA", "NA", "NA", "NA")
Set-RowValues $ws 30 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws 31 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws 32 @(71000, -26000, -32000, -74000, -51000, -2000, 31000, 6000, -237800, -29700)
Set-RowValues $ws 33 @(67000, 65000, -18000, 357000, 127000, -289000, -285000, 54000, 109700, -100400)
Set-RowValues $ws 34 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws 35 @(67000, 65000, -18000, 357000, 127000, -289000, -285000, 54000, 109700, -100400)
Set-RowValues $ws 38 @(43465, 43373, 43281, 43190, 43100, 43008, 42916, 42825, 42735, 42643)
Set-RowValues $ws 41 @(981000, 989000, 874000, 715000, 722000, 919000, 796000, 923000, 875800, 990100)
Set-RowValues $ws 42 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws 43 @(585000, 246000, 280000, 608000, 371000, 265000, 284000, 290000, 217900, 154200)
Set-RowValues $ws 44 @(316000, 298000, 233000, 123000, 243000, 133000, 150000, 113000, 160200, 63900)
Set-RowValues $ws 45 @(2352000, 2137000, 2579000, 1822000, 2033000, 1714000, 1138000, 1093000, 984400, 896600)
Set-RowValues $ws 46 @(4234000, 3670000, 3966000, 3268000, 3369000, 3031000, 2368000, 2419000, 2238300, 2104700)
Set-RowValues $ws 47 @(115000, 113000, 90000, 92000, 93000, 104000, 101000, 38000, 41300, 42300)
Set-RowValues $ws 48 @(27245000, 26499000, 25760000, 24474000, 23978000, 23466000, 22904000, 22016000, 20635300, 19891700)
Set-RowValues $ws 49 @(77000, 77000, 77000, 77000, 77000, 77000, 77000, 77000, 76800, 76800)
Set-RowValues $ws 50 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws 51 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws 52 @(316000, 381000, 434000, 430000, 389000, 456000, 1150000, 1506000, 711000, 573600)
Set-RowValues $ws 53 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws 54 @(31987000, 30740000, 30327000, 28341000, 27906000, 27134000, 26600000, 26056000, 23702700, 22689100)
Set-RowValues $ws 57 @(58000, 80000, 29000, 21000, 25000, 59000, 62000, 54000, 48600, 38600)
Set-RowValues $ws 58 @("NA", 66000, 137000, 0, 0, 41000, "NA", "NA", 247500, "NA")
Set-RowValues $ws 59 @(1684000, 1204000, 1562000, 874000, 1226000, 911000, 778000, 817000, 780600, 2570300)
Set-RowValues $ws 60 @(1742000, 1350000, 1728000, 895000, 1251000, 1011000, 840000, 871000, 1076700, 2608900)
Set-RowValues $ws 61 @(28236000, 27467000, 26782000, 25656000, 25337000, 24925000, 24657000, 24088000, 21687500, 19033500)
Set-RowValues $ws 62 @(80000, 92000, 83000, 83000, 78000, 115000, 99000, 102000, 100100, 341000)
Set-RowValues $ws 63 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws 64 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws 65 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws 66 @(32513000, 31349000, 31344000, 29738000, 29670000, 29027000, 28213000, 27392000, 25099200, 24291800)
Set-RowValues $ws 68 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws 69 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws 70 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws 71 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws 72 @(-4156000, -4223000, -4288000, -4270000, -4627000, -4754000, -4465000, -4180000, -4233900, -4343600)
Set-RowValues $ws 73 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws 74 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws 75 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws 76 @(-526000, -609000, -1017000, -1397000, -1764000, -1893000, -1613000, -1336000, -1396400, -1602700)
Set-RowValues $ws 77 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws 80 @(43465, 43373, 43281, 43190, 43100, 43008, 42916, 42825, 42735, 42643)
Set-RowValues $ws 81 @(67000, 65000, -18000, 357000, 127000, -289000, -285000, 54000, 109700, -100400)
Set-RowValues $ws 83 @(116000, 113000, 111000, 109000, 104000, 92000, 90000, 70000, 68000, 49200)
Set-RowValues $ws 84 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws 85 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws 86 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws 87 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws 88 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws 89 @(486000, 522000, 513000, 469000, 336000, 359000, 227000, 309000, -204800, -55000)
Set-RowValues $ws 91 @(-931000, -1204000, -732000, -776000, -454000, -565000, -1019000, -1319000, 2542600, -4620700)
Set-RowValues $ws 92 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws 93 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws 94 @(-932000, -1230000, -716000, -776000, -455000, -569000, -1067000, -1290000, -4401200, 7100)
Set-RowValues $ws 96 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws 97 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws 98 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws 99 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws 100 @(670000, 369000, 1052000, 116000, 157000, 241000, 443000, 2095000, 976600, 3919300)
Set-RowValues $ws 101 @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws 102 @(224000, -339000, 849000, -191000, 38000, 31000, -397000, 1114000, -22700, 264600)

# A handful of shifted cells carry slightly restated figures versus a pure
# shift of the prior quarter's value (minor vendor restatements).
$ws.Range("H61").Value = 25337000
$ws.Range("I61").Value = 24925000
$ws.Range("J61").Value = 24657000
$ws.Range("H62").Value = 78000
$ws.Range("I62").Value = 115000
$ws.Range("J62").Value = 99000
